$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dq``"
$ws.Range("B2").Value = "fdawd"
$ws.Range("C2").Value = "fqw"
$ws.Range("D2").Value = "fqw"
